# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Tue Sep 24 07:34:26 UTC 2024 with GitHub Actions"
# Price (col D) and 1h volume-change (col E) cells are updated in place,
# and the Bittensor/Kaspa rows (28/29) swap order.
# A leading apostrophe forces each assignment to stay a Text cell
# (matching the original inlineStr cells) instead of Excel auto-detecting
# a number for values like "1.00" or "601.49".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.427.62"
$ws.Range('E2').Value = "'  -0.05%  "
$ws.Range('D3').Value = "'2.642.07"
$ws.Range('E3').Value = "'  -0.03%  "
$ws.Range('E4').Value = "'  +0.07%  "
$ws.Range('D5').Value = "'601.49"
$ws.Range('E5').Value = "'  +1.81%  "
$ws.Range('D6').Value = "'146.30"
$ws.Range('E6').Value = "'  +0.23%  "
$ws.Range('E7').Value = "'  +0.05%  "
$ws.Range('E8').Value = "'  -0.53%  "
$ws.Range('E9').Value = "'  +0.98%  "
$ws.Range('E10').Value = "'  -0.76%  "
$ws.Range('D11').Value = "'0.367"
$ws.Range('E11').Value = "'  +3.40%  "
$ws.Range('E12').Value = "'  +0.15%  "
$ws.Range('D13').Value = "'27.35"
$ws.Range('E13').Value = "'  -1.06%  "
$ws.Range('D14').Value = "'3.123.70"
$ws.Range('E14').Value = "'  +0.25%  "
$ws.Range('D15').Value = "'63.307.08"
$ws.Range('E15').Value = "'  -0.09%  "
$ws.Range('E16').Value = "'  -0.77%  "
$ws.Range('D17').Value = "'2.630.08"
$ws.Range('E17').Value = "'  +0.11%  "
$ws.Range('D18').Value = "'11.45"
$ws.Range('E18').Value = "'  +0.94%  "
$ws.Range('D19').Value = "'4.52"
$ws.Range('E19').Value = "'  +3.66%  "
$ws.Range('D20').Value = "'341.46"
$ws.Range('E20').Value = "'  -0.03%  "
$ws.Range('E21').Value = "'  +2.82%  "
$ws.Range('E22').Value = "'  +0.04%  "
$ws.Range('E23').Value = "'  -3.42%  "
$ws.Range('D24').Value = "'66.78"
$ws.Range('E24').Value = "'  -1.67%  "
$ws.Range('D25').Value = "'1.70"
$ws.Range('E25').Value = "'  +0.88%  "
$ws.Range('E26').Value = "'  +4.36%  "
$ws.Range('D27').Value = "'1.54"
$ws.Range('E27').Value = "'  -2.25%  "
$ws.Range('B28').Value = "'Kaspa"
$ws.Range('C28').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('D28').Value = "'0.163"
$ws.Range('E28').Value = "'  -1.80%  "
$ws.Range('B29').Value = "'Bittensor"
$ws.Range('C29').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('D29').Value = "'547.15"
$ws.Range('E29').Value = "'  -1.42%  "
$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = "'  -0.26%  "
$ws.Range('D31').Value = "'7.87"
$ws.Range('E31').Value = "'  +0.42%  "
$ws.Range('D32').Value = "'2.05"
$ws.Range('E32').Value = "'  +4.61%  "
$ws.Range('E33').Value = "'  -2.77%  "
$ws.Range('D34').Value = "'0.0₃0805"
$ws.Range('E34').Value = "'  -0.80%  "
$ws.Range('D35').Value = "'5.23"
$ws.Range('E35').Value = "'  +6.49%  "
$ws.Range('D36').Value = "'168.11"
$ws.Range('E36').Value = "'  -4.12%  "
$ws.Range('D37').Value = "'0.405"
$ws.Range('E37').Value = "'  +0.78%  "
$ws.Range('E38').Value = "'  -0.01%  "
$ws.Range('D39').Value = "'19.08"
$ws.Range('E39').Value = "'  -0.43%  "
$ws.Range('E40').Value = "'  +6.47%  "
$ws.Range('E41').Value = "'  -0.05%  "
$ws.Range('D42').Value = "'168.95"
$ws.Range('E42').Value = "'  -1.09%  "
$ws.Range('E43').Value = "'  +0.47%  "
$ws.Range('D44').Value = "'22.49"
$ws.Range('E44').Value = "'  +2.74%  "
$ws.Range('E45').Value = "'  +4.18%  "
$ws.Range('D46').Value = "'0.626"
$ws.Range('E46').Value = "'  -0.45%  "
$ws.Range('D47').Value = "'0.0246"
$ws.Range('E47').Value = "'  +2.62%  "
$ws.Range('D48').Value = "'0.0962"
$ws.Range('E48').Value = "'  +0.29%  "
$ws.Range('D49').Value = "'18.86"
$ws.Range('E49').Value = "'  +0.27%  "
$ws.Range('D50').Value = "'1.79"
$ws.Range('E50').Value = "'  +4.73%  "
$ws.Range('D51').Value = "'11.27"
$ws.Range('E51').Value = "'  -0.59%  "

Write-Output "Applied 87 cell edits"
